# "Regenerate graphs with MBE"
#
# The reviewer comment recommending MAE over RMSE is struck through, and a
# new bold author's-response paragraph is added right after it, explaining
# that MBE (rather than MAE) was used so the results could be compared
# against Best & Grimmond (2012).

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "I suggest calculation of MAE instead of RMSE.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the MAE/RMSE reviewer comment paragraph"
}

$commentPara = $rng.Paragraphs(1)

# Add the new (bold) author's-response paragraph directly after the comment,
# before striking the comment through, so the new paragraph does not inherit
# the strike-through formatting.
$commentPara.Range.InsertParagraphAfter()
$responsePara = $commentPara.Next()
$responsePara.Range.Text = "Added calculations of MBE (instead of MAE) in order to compare to Best & Grimmond (2012) who uses RMSE and MBE."
$responsePara.Range.Font.Bold = $true
$responsePara.Range.Font.BoldBi = $true

# Strike through the original reviewer comment (applies to the run and to
# the paragraph mark).
$commentPara.Range.Font.StrikeThrough = $true
